$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 00:25"

# Update numeric stats for countries with new data
$ws.Cells.Item(4, 2).Value = 274987
$ws.Cells.Item(4, 3).Value = 30110
$ws.Cells.Item(4, 5).Value = 255878
$ws.Cells.Item(4, 7).Value = 994
$ws.Cells.Item(4, 8).Value = 7065
$ws.Cells.Item(17, 2).Value = 11524
$ws.Cells.Item(17, 3).Value = 395
$ws.Cells.Item(17, 5).Value = 9334
$ws.Cells.Item(77, 6).Value = 21
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 8).Value = 2
$ws.Cells.Item(165, 3).Value = 3
$ws.Cells.Item(165, 4).Value = 2
$ws.Cells.Item(165, 8).Value = 0
$ws.Cells.Item(167, 2).Value = 15
$ws.Cells.Item(167, 3).Value = 6
$ws.Cells.Item(167, 4).Value = 0
$ws.Cells.Item(167, 5).Value = 15
$ws.Cells.Item(167, 6).Value = 1
$ws.Cells.Item(168, 4).Value = 2
$ws.Cells.Item(168, 5).Value = 12
$ws.Cells.Item(169, 2).Value = 14
$ws.Cells.Item(169, 4).Value = 3
$ws.Cells.Item(169, 5).Value = 11
$ws.Cells.Item(170, 2).Value = 13
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 1
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(172, 2).Value = 12
$ws.Cells.Item(172, 3).Value = 2
$ws.Cells.Item(172, 5).Value = 12
$ws.Cells.Item(172, 6).Value = 2
$ws.Cells.Item(172, 8).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 10
$ws.Cells.Item(174, 2).Value = 11
$ws.Cells.Item(174, 4).Value = 3
$ws.Cells.Item(174, 5).Value = 7
$ws.Cells.Item(174, 8).Value = 1
$ws.Cells.Item(177, 5).Value = 10
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 5).Value = 9
$ws.Cells.Item(178, 7).Value = 1
$ws.Cells.Item(178, 8).Value = 1
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 3
$ws.Cells.Item(179, 5).Value = 7
$ws.Cells.Item(179, 8).Value = 0
$ws.Cells.Item(180, 2).Value = 10
$ws.Cells.Item(180, 3).Value = 2
$ws.Cells.Item(180, 4).Value = 2
$ws.Cells.Item(180, 5).Value = 6
$ws.Cells.Item(180, 8).Value = 2
$ws.Cells.Item(191, 3).Value = 2
$ws.Cells.Item(192, 2).Value = 6
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 8).Value = 1
$ws.Cells.Item(193, 4).Value = 1
$ws.Cells.Item(193, 8).Value = 0
$ws.Cells.Item(195, 4).Value = 2
$ws.Cells.Item(195, 5).Value = 3
$ws.Cells.Item(195, 8).Value = 1
$ws.Cells.Item(196, 2).Value = 5
$ws.Cells.Item(196, 4).Value = 0
$ws.Cells.Item(196, 5).Value = 5
$ws.Cells.Item(196, 8).Value = 0
$ws.Cells.Item(197, 5).Value = 4
$ws.Cells.Item(197, 8).Value = 1
$ws.Cells.Item(198, 4).Value = 2
$ws.Cells.Item(198, 5).Value = 3
$ws.Cells.Item(198, 8).Value = 0
$ws.Cells.Item(199, 2).Value = 4
$ws.Cells.Item(199, 3).Value = 1
$ws.Cells.Item(199, 4).Value = 0
$ws.Cells.Item(199, 5).Value = 4
$ws.Cells.Item(200, 3).Value = 0
$ws.Cells.Item(200, 5).Value = 3
$ws.Cells.Item(200, 8).Value = 1
$ws.Cells.Item(201, 4).Value = 2
$ws.Cells.Item(201, 5).Value = 1
$ws.Cells.Item(202, 2).Value = 3
$ws.Cells.Item(202, 4).Value = 0
$ws.Cells.Item(202, 5).Value = 3
$ws.Cells.Item(202, 8).Value = 0
$ws.Cells.Item(206, 3).Value = 1
$ws.Cells.Item(206, 4).Value = 1
$ws.Cells.Item(206, 5).Value = 2

# Re-sort causes these rows to now reference different countries;
# update country names (column A) to reflect the new sort order
$ws.Cells.Item(164, 1).Value = "Siria"
$ws.Cells.Item(165, 1).Value = "Benin"
$ws.Cells.Item(167, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(168, 1).Value = "Mongolia"
$ws.Cells.Item(169, 1).Value = "Namibia"
$ws.Cells.Item(170, 1).Value = "Santa Lucia"
$ws.Cells.Item(172, 1).Value = "Granada"
$ws.Cells.Item(173, 1).Value = "Libia"
$ws.Cells.Item(174, 1).Value = "Curazao"
$ws.Cells.Item(176, 1).Value = "Mozambique"
$ws.Cells.Item(177, 1).Value = "Seychelles"
$ws.Cells.Item(178, 1).Value = "Surinam"
$ws.Cells.Item(179, 1).Value = "Groenlandia"
$ws.Cells.Item(180, 1).Value = "Sudan"
$ws.Cells.Item(188, 1).Value = "Fiyi"
$ws.Cells.Item(190, 1).Value = "Santa Sede"
$ws.Cells.Item(191, 1).Value = "Somalia"
$ws.Cells.Item(192, 1).Value = "Cabo Verde"
$ws.Cells.Item(193, 1).Value = "Nepal"
$ws.Cells.Item(194, 1).Value = "San Bartolome"
$ws.Cells.Item(195, 1).Value = "Mauritania"
$ws.Cells.Item(196, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(197, 1).Value = "Nicaragua"
$ws.Cells.Item(198, 1).Value = "Butan"
$ws.Cells.Item(199, 1).Value = "Belice"
$ws.Cells.Item(200, 1).Value = "Botsuana"
$ws.Cells.Item(201, 1).Value = "Gambia"
$ws.Cells.Item(202, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(205, 1).Value = "Burundi"
$ws.Cells.Item(206, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(207, 1).Value = "Sierra Leona"
$ws.Cells.Item(208, 1).Value = "Bonaire, San Eustaquio y Saba"
